$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows above the old row 2 (pushes old rows 2-30 down to 4-32)
$ws.Range("A2:N3").EntireRow.Insert()

# The inserted rows pick up formatting copied from the row above (row 1).
# Remove that formatting so the new rows 2 & 3 are plain/unstyled, matching
# the rest of the data rows.
$ws.Range("A2:N3").ClearFormats()

# Row 1: replace the text headers with numeric column-index codes (0-13).
# Keep the existing bold/border/centered style (style index 1) which stays
# attached to row 1.
$ws.Range("A1").Value = 0
$ws.Range("B1").Value = 1
$ws.Range("C1").Value = 2
$ws.Range("D1").Value = 3
$ws.Range("E1").Value = 4
$ws.Range("F1").Value = 5
$ws.Range("G1").Value = 6
$ws.Range("H1").Value = 7
$ws.Range("I1").Value = 8
$ws.Range("J1").Value = 9
$ws.Range("K1").Value = 10
$ws.Range("L1").Value = 11
$ws.Range("M1").Value = 12
$ws.Range("N1").Value = 13

# Row 2 (new): only E2 gets a value.
$ws.Range("E2").Value = "Washer"

# Row 3 (new): re-create the original header row text labels (the old
# row 1), minus the thread_size / material_surface / Pkg. Qty. sub-code
# labels that are dropped (K3, M3, N3 stay blank).
$ws.Range("A3").Value = "Lg."
$ws.Range("B3").Value = "Threading"
$ws.Range("C3").Value = "HeadDia."
$ws.Range("D3").Value = "Head Ht."
$ws.Range("E3").Value = "OD"
$ws.Range("F3").Value = "Thick."
$ws.Range("G3").Value = "DriveSize"
$ws.Range("H3").Value = "TensileStrength, psi"
$ws.Range("I3").Value = "Specifications Met"
$ws.Range("J3").Value = "Pkg.Qty."
$ws.Range("L3").Value = "Pkg."
